$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H32").Value = 832.7778
$ws.Range("I32").Value = 649.5
$ws.Range("K32").Value = 649.5
$ws.Range("M32").Value = -323.5

$ws.Range("H131").Value = 3807.3
$ws.Range("I131").Value = 3786
$ws.Range("K131").Value = 11358
$ws.Range("M131").Value = -6318

$ws.Range("H135").Value = 3226553.2
$ws.Range("I135").Value = 4000668.2
$ws.Range("K135").Value = 36006013.8
$ws.Range("M135").Value = -36003478.8

$ws.Range("H137").Value = 5820569.5
$ws.Range("I137").Value = 11366550
$ws.Range("K137").Value = 34099650
$ws.Range("M137").Value = -34097100

$ws.Range("H138").Value = 2740.0222
$ws.Range("I138").Value = 1952.6522
$ws.Range("J138").Value = 3563.182
$ws.Range("K138").Value = 5857.9566
$ws.Range("L138").Value = 10689.546
$ws.Range("M138").Value = -717.9565999999995
$ws.Range("N138").Value = -20969.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 6949.5
$ws.Range("I41").Value = 6949.5
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 6949.5
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -6535.5
$ws.Range("N41").ClearContents()

$ws.Range("H45").Value = 2500
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H74").Value = 5262.0713
$ws.Range("I74").Value = 2298.3333
$ws.Range("J74").Value = 6070.364
$ws.Range("K74").Value = 2298.3333
$ws.Range("L74").Value = 6070.364
$ws.Range("M74").Value = -1424.3333
$ws.Range("N74").Value = -7818.364

$ws.Range("H77").Value = 5262.0713
$ws.Range("I77").Value = 2298.3333
$ws.Range("J77").Value = 6070.364
$ws.Range("K77").Value = 11491.6665
$ws.Range("L77").Value = 30351.82
$ws.Range("M77").Value = -7123.666499999999
$ws.Range("N77").Value = -39087.82

$ws.Range("H92").Value = 10000000
$ws.Range("J92").Value = 10000000
$ws.Range("L92").Value = 10000000
$ws.Range("N92").Value = -10004992

$ws.Range("H112").Value = 37899.2
$ws.Range("J112").Value = 37899.2
$ws.Range("L112").Value = 37899.2
$ws.Range("N112").Value = -40853.2

$ws.Range("H132").Value = 6460.7393
$ws.Range("I132").Value = 4347.2104
$ws.Range("K132").Value = 13041.6312
$ws.Range("M132").Value = -10511.6312

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H94").Value = 790.71875
$ws.Range("I94").Value = 703.96155
$ws.Range("J94").Value = 1166.6666
$ws.Range("K94").Value = 703.96155
$ws.Range("L94").Value = 1166.6666
$ws.Range("M94").Value = -252.96155
$ws.Range("N94").Value = -2068.6666

$ws.Range("H134").Value = 7338.1
$ws.Range("J134").Value = 8071.8
$ws.Range("L134").Value = 24215.4
$ws.Range("N134").Value = -29285.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3143.2856
$ws.Range("I16").Value = 3333.3333
$ws.Range("J16").Value = 3000.75
$ws.Range("K16").Value = 3333.3333
$ws.Range("L16").Value = 3000.75
$ws.Range("M16").Value = -3046.3333
$ws.Range("N16").Value = -3574.75

$ws.Range("H31").Value = 3442.898
$ws.Range("I31").Value = 1806
$ws.Range("K31").Value = 1806
$ws.Range("M31").Value = -1511

$ws.Range("H34").Value = 3442.898
$ws.Range("I34").Value = 1806
$ws.Range("K34").Value = 1806
$ws.Range("M34").Value = -1604

$ws.Range("H58").Value = 4437.9375
$ws.Range("I58").Value = 3165.4285
$ws.Range("K58").Value = 3165.4285
$ws.Range("M58").Value = -2962.4285

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H105").Value = 33684.332
$ws.Range("I105").Value = 33684.332
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 33684.332
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -31937.332
$ws.Range("N105").ClearContents()

$ws.Range("H113").Value = 3143.2856
$ws.Range("I113").Value = 3333.3333
$ws.Range("J113").Value = 3000.75
$ws.Range("K113").Value = 3333.3333
$ws.Range("L113").Value = 3000.75
$ws.Range("M113").Value = -1163.3333
$ws.Range("N113").Value = -7340.75

$ws.Range("H132").Value = 50558.703
$ws.Range("I132").Value = 3242.1428
$ws.Range("K132").Value = 9726.428400000001
$ws.Range("M132").Value = -7196.428400000001

$ws.Range("H136").Value = 4437.9375
$ws.Range("I136").Value = 3165.4285
$ws.Range("K136").Value = 9496.2855
$ws.Range("M136").Value = -6946.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2230.6
$ws.Range("I34").Value = 324.66666
$ws.Range("J34").Value = 3501.2222
$ws.Range("K34").Value = 973.9999799999999
$ws.Range("L34").Value = 10503.6666
$ws.Range("M34").Value = -889.9999799999999
$ws.Range("N34").Value = -10671.6666

$ws.Range("H39").Value = 7211.5557
$ws.Range("J39").Value = 7488
$ws.Range("L39").Value = 22464
$ws.Range("N39").Value = -23052

$ws.Range("H69").Value = 1000
$ws.Range("I69").Value = 1000
$ws.Range("K69").Value = 3000
$ws.Range("M69").Value = -2189

$ws.Range("H72").Value = 1000
$ws.Range("I72").Value = 1000
$ws.Range("K72").Value = 9000
$ws.Range("M72").Value = -4944

$ws.Range("H75").Value = 794.1818
$ws.Range("J75").Value = 738.5714
$ws.Range("L75").Value = 2215.7142
$ws.Range("N75").Value = -4211.7142

$ws.Range("H78").Value = 794.1818
$ws.Range("J78").Value = 738.5714
$ws.Range("L78").Value = 6647.1426
$ws.Range("N78").Value = -16631.1426

$ws.Range("H103").Value = 1160.4286
$ws.Range("I103").Value = 1278.8334
$ws.Range("K103").Value = 3836.5002
$ws.Range("M103").Value = -2957.5002

$ws.Range("H129").Value = 2397.889
$ws.Range("J129").Value = 2798
$ws.Range("L129").Value = 8394
$ws.Range("N129").Value = -18394

$ws.Range("H131").Value = 4023.4194
$ws.Range("I131").Value = 980.3333
$ws.Range("K131").Value = 2940.9999
$ws.Range("M131").Value = 2099.0001

$ws.Range("H137").Value = 1694.3125
$ws.Range("I137").Value = 1181
$ws.Range("J137").Value = 2207.625
$ws.Range("K137").Value = 3543
$ws.Range("L137").Value = 6622.875
$ws.Range("M137").Value = 1557
$ws.Range("N137").Value = -16822.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 22706.79
$ws.Range("I24").Value = 12221
$ws.Range("J24").Value = 27546.385
$ws.Range("K24").Value = 12221
$ws.Range("L24").Value = 27546.385
$ws.Range("M24").Value = -12048
$ws.Range("N24").Value = -27892.385

$ws.Range("H96").Value = 52489.5
$ws.Range("J96").Value = 52489.5
$ws.Range("L96").Value = 52489.5
$ws.Range("N96").Value = -57981.5

$ws.Range("H132").Value = 6185.536
$ws.Range("I132").Value = 4039.8655
$ws.Range("K132").Value = 12119.5965
$ws.Range("M132").Value = -9589.5965

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3329.1667
$ws.Range("I16").Value = 3205.2
$ws.Range("K16").Value = 3205.2
$ws.Range("M16").Value = -3035.2

$ws.Range("H58").Value = 6910
$ws.Range("J58").Value = 7949.6665
$ws.Range("L58").Value = 7949.6665
$ws.Range("N58").Value = -8469.666499999999

$ws.Range("H100").Value = 7816475.5
$ws.Range("I100").Value = 14709339
$ws.Range("K100").Value = 14709339
$ws.Range("M100").Value = -14708798

$ws.Range("H132").Value = 5160.88
$ws.Range("I132").Value = 3901.1
$ws.Range("J132").Value = 10200
$ws.Range("K132").Value = 11703.3
$ws.Range("L132").Value = 30600
$ws.Range("M132").Value = -9173.299999999999
$ws.Range("N132").Value = -35660

$ws.Range("H136").Value = 5722.3687
$ws.Range("J136").Value = 8452
$ws.Range("L136").Value = 25356
$ws.Range("N136").Value = -30456

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9627.111000000001
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 9627.111000000001
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 9627.111000000001
$ws.Range("N41").Value = -10407.111
$ws.Range("M41").ClearContents()

$ws.Range("H61").Value = 3805.3076
$ws.Range("I61").Value = 3980.75
$ws.Range("J61").Value = 1700
$ws.Range("K61").Value = 3980.75
$ws.Range("L61").Value = 1700
$ws.Range("M61").Value = -3688.75
$ws.Range("N61").Value = -2284

$ws.Range("H81").Value = 2333.7222
$ws.Range("I81").Value = 1615.7333
$ws.Range("J81").Value = 5923.6665
$ws.Range("K81").Value = 3231.4666
$ws.Range("L81").Value = 11847.333
$ws.Range("M81").Value = -2170.4666
$ws.Range("N81").Value = -13969.333

$ws.Range("H84").Value = 2333.7222
$ws.Range("I84").Value = 1615.7333
$ws.Range("J84").Value = 5923.6665
$ws.Range("K84").Value = 16157.333
$ws.Range("L84").Value = 59236.665
$ws.Range("M84").Value = -10853.333
$ws.Range("N84").Value = -69844.66500000001

$ws.Range("H136").Value = 7446.8
$ws.Range("I136").Value = 6710.778
$ws.Range("J136").Value = 8550.833000000001
$ws.Range("K136").Value = 20132.334
$ws.Range("L136").Value = 25652.499
$ws.Range("M136").Value = -17582.334
$ws.Range("N136").Value = -30752.499
